$d = $word.ActiveDocument

# --- Step 1: move the _GoBack bookmark from paragraph 1 to paragraph 4 (its
#     eventual home is the now-empty list paragraph near the end of the list).
$p4 = $d.Paragraphs.Item(4).Range
$d.Bookmarks.Add("_GoBack", $p4)

# --- Step 2: text edits on the (still) first paragraph.
$d.Content.Find.Execute("Do I need to get rid of common ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "How to have bar graphs show only the score value as ", 2)
$d.Content.Find.Execute("words", $true, $false, $false, $false, $false,
                         $true, 1, $false, "category", 2)

# --- Step 3: append a bare trailing paragraph after paragraph 4 (currently
#     the last paragraph in the document / holds the relocated bookmark).
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$pLast.InsertParagraphAfter()
$newTail = $d.Paragraphs.Item($d.Paragraphs.Count)
$newTail.Range.ListFormat.RemoveNumbers()
$newTail.Style = "Normal"

# --- Step 4: insert the new list item "Do I need to get rid of common words"
#     right before the (old) second paragraph ("What is problem with ...").
$p2 = $d.Paragraphs.Item(2).Range
$p2.InsertParagraphBefore()
$newListPara = $d.Paragraphs.Item(2)
$newListPara.Range.Text = "Do I need to get rid of common words"

# --- Step 5: insert two bare blank paragraphs between paragraph 1 and the
#     new list item created above.
$p1 = $d.Paragraphs.Item(1).Range
$p1.InsertParagraphAfter()
$blank1 = $d.Paragraphs.Item(2)
$blank1.Range.ListFormat.RemoveNumbers()
$blank1.Style = "Normal"

$blank1.Range.InsertParagraphAfter()
$blank2 = $d.Paragraphs.Item(3)
$blank2.Range.ListFormat.RemoveNumbers()
$blank2.Style = "Normal"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i : [" $p.Range.Text "]"
}
